# Script output update: the "Numéro de page" value that used to sit on
# 2025-03-11 (row 9) didn't apply to that row anymore, and a new day
# (2025-03-12) was appended with that same "NA" value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 no longer has a page number -> clear it to an empty string, the
# same "empty" shape already used by rows 2-8's column C (a lone leading
# apostrophe forces an explicit, literal empty-text entry instead of
# fully blanking the cell).
$ws.Range("C9").Value = "'"

# New row 10.
# A10 looks like a date ("2025-03-12"); prefix with an apostrophe so Excel
# stores it as literal text, matching the other Date-column cells, instead
# of auto-converting it to a date serial number.
$ws.Range("A10").Value = "'2025-03-12"
$ws.Range("B10").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C10").Value = "NA"
$ws.Range("D10").Value = 1
